# Insert a new weekly data row for "Apio" (Macroferia Regional de Talca) at row 296.
# This pushes the existing rows 296-337 down to 297-338 and extends the sheet
# dimension from A1:R337 to A1:R338.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row before the current row 296 - everything below shifts down.
$ws.Rows.Item(296).Insert()

# Populate the newly inserted row 296 with the new market observation.
$ws.Cells.Item(296, 1).Value  = 5
$ws.Cells.Item(296, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(296, 3).Value  = "Maule"
$ws.Cells.Item(296, 4).Value  = 45131
$ws.Cells.Item(296, 5).Value  = 7
$ws.Cells.Item(296, 6).Value  = 100112017
$ws.Cells.Item(296, 7).Value  = "Apio"
$ws.Cells.Item(296, 8).Value  = "Americana (o)"
$ws.Cells.Item(296, 9).Value  = "Primera"
$ws.Cells.Item(296, 10).Value = 700
$ws.Cells.Item(296, 11).Value = 5500
$ws.Cells.Item(296, 12).Value = 5500
$ws.Cells.Item(296, 13).Value = 5500
$ws.Cells.Item(296, 14).Value = "$/docena de matas"
$ws.Cells.Item(296, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(296, 16).Value = 917
$ws.Cells.Item(296, 17).Value = 6
$ws.Cells.Item(296, 18).Value = "Hortaliza"
